# Add a new row (row 53) of portfolio data for 2025-10-07 to the bottom
# of the existing table on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as plain text (matching the existing rows, which
# are stored as inline/shared strings rather than real Excel dates). Simply
# assigning a date-looking string via .Value triggers Excel's automatic
# date recognition, turning it into a date serial number with date
# formatting. To avoid that, build the text through a formula (a string
# literal, so it is never treated as a date) and then convert that
# computed value into a static value with Copy + PasteSpecial (values
# only), which keeps it as plain text without picking up a new/different
# cell style.
$dateCell = $ws.Range("A53")
$dateCell.Formula = "=""2025-10-07"""
$dateCell.Copy()
$dateCell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# Numeric columns can be set directly.
$ws.Range("B53").Value = 54.0099983215332
$ws.Range("C53").Value = 698.0499877929688
$ws.Range("D53").Value = 337.8500061035156
